$wb = $excel.ActiveWorkbook

# Insert the new "ExistingRunConfig_Data" sheet right after "DescriptiveStatistics_Data"
# (and therefore right before "API_Data"), matching the new <sheets> order in workbook.xml.
$afterSheet = $wb.Worksheets.Item("DescriptiveStatistics_Data")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "ExistingRunConfig_Data"

# Column A values (Level labels). Row 10 and row 11 both use "Level 10".
$levels = @("Level 1","Level 2","Level 3","Level 4","Level 5","Level 6","Level 7","Level 8","Level 9","Level 10","Level 10","Level 11","Level 12","Level 13","Level 14","Level 15","Level 16","Level 17","Level 18","Level 19")

# Column B/C/D/E numeric values, one tuple per row.
$nums = @(
    @(1,1,1,1),
    @(1,1,1,1),
    @(1,1,1,1),
    @(3,4,2,2),
    @(1,24,1,1),
    @(2,2,2,2),
    @(1,7,5,1),
    @(1,3,2,4),
    @(1,2,1,1),
    @(99,19,19,19),
    @(99,19,19,7),
    @(2,5,2,2),
    @(3,2,3,2),
    @(11,11,11,11),
    @(11,11,11,11),
    @(1,1,1,1),
    @(3,7,5,7),
    @(1,8,5,1),
    @(1,1,1,1),
    @(2,1,1,1)
)

# Write rows 1-12 first (A + default F="NO"), then flip row 1's F to "YES",
# then write the remaining rows 13-20 -- this reproduces the exact shared-string
# insertion order (Level 1, NO, Level 2 ... Level 11, YES, Level 12 ... Level 19).
for ($i = 0; $i -lt 12; $i++) {
    $r = $i + 1
    $newSheet.Cells.Item($r, 1).Value = $levels[$i]
    $newSheet.Cells.Item($r, 2).Value = $nums[$i][0]
    $newSheet.Cells.Item($r, 3).Value = $nums[$i][1]
    $newSheet.Cells.Item($r, 4).Value = $nums[$i][2]
    $newSheet.Cells.Item($r, 5).Value = $nums[$i][3]
    $newSheet.Cells.Item($r, 6).Value = "NO"
}

$newSheet.Cells.Item(1, 6).Value = "YES"

for ($i = 12; $i -lt 20; $i++) {
    $r = $i + 1
    $newSheet.Cells.Item($r, 1).Value = $levels[$i]
    $newSheet.Cells.Item($r, 2).Value = $nums[$i][0]
    $newSheet.Cells.Item($r, 3).Value = $nums[$i][1]
    $newSheet.Cells.Item($r, 4).Value = $nums[$i][2]
    $newSheet.Cells.Item($r, 5).Value = $nums[$i][3]
    $newSheet.Cells.Item($r, 6).Value = "NO"
}
